$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.453.84'
$ws.Range("E2").Value = '  +5.43%  '
$ws.Range("D3").Value = '3.777.08'
$ws.Range("E3").Value = '  +22.60%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '''617.77'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +8.35%  '
$ws.Range("D6").Value = '''177.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.37%  '
$ws.Range("D7").Value = '3.778.27'
$ws.Range("E7").Value = '  +22.65%  '
$ws.Range("E9").Value = '  +6.73%  '
$ws.Range("D10").Value = '''0.170'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +12.93%  '
$ws.Range("D11").Value = '''6.42'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.47%  '
$ws.Range("E12").Value = '  +8.72%  '
$ws.Range("D13").Value = '''40.59'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.51%  '
$ws.Range("D14").Value = '''0.0000261'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.14%  '
$ws.Range("D15").Value = '4.409.13'
$ws.Range("E15").Value = '  +22.68%  '
$ws.Range("D16").Value = '3.792.57'
$ws.Range("E16").Value = '  +23.07%  '
$ws.Range("D17").Value = '70.587.05'
$ws.Range("E17").Value = '  +5.70%  '
$ws.Range("E18").Value = '  +1.52%  '
$ws.Range("E19").Value = '  +9.21%  '
$ws.Range("D20").Value = '''525.13'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +8.72%  '
$ws.Range("D21").Value = '''16.98'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.51%  '
$ws.Range("D22").Value = '''9.41'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +22.83%  '
$ws.Range("D23").Value = '''0.747'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.39%  '
$ws.Range("D24").Value = '''88.08'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +5.81%  '
$ws.Range("E25").Value = '  +11.81%  '
$ws.Range("E26").Value = '  +7.49%  '
$ws.Range("D27").Value = '''10.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.21%  '
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '''0.0000123'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +33.32%  '
$ws.Range("E30").Value = '  +10.10%  '
$ws.Range("D31").Value = '''2.91'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +12.88%  '
$ws.Range("E32").Value = '  +2.50%  '
$ws.Range("D33").Value = '''32.19'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +15.59%  '
$ws.Range("E34").Value = '  +4.26%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("E36").Value = '  +11.98%  '
$ws.Range("D37").Value = '''1.04'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +10.97%  '
$ws.Range("D38").Value = '''0.343'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +10.49%  '
$ws.Range("D39").Value = '''2.20'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.47%  '
$ws.Range("E40").Value = '  +8.81%  '
$ws.Range("D41").Value = '''51.57'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.45%  '
$ws.Range("B42").Value = 'Bittensor'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D42").Value = '''430.15'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +17.26%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.153.23'
$ws.Range("E43").Value = '  +13.24%  '
$ws.Range("D44").Value = '''8.88'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = '''44.59'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.41%  '
$ws.Range("E46").Value = '  +4.34%  '
$ws.Range("E47").Value = '  +8.24%  '
$ws.Range("D48").Value = '''27.79'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +8.52%  '
$ws.Range("D49").Value = '''140.70'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.23%  '
$ws.Range("E50").Value = '  +11.75%  '
$ws.Range("E51").Value = '  +0.02%  '
